$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(32, 8).Value = 3354.147
$ws.Cells.Item(32, 10).Value = 3543.4688
$ws.Cells.Item(32, 12).Value = 3543.4688
$ws.Cells.Item(32, 14).Value = -4195.468800000001
$ws.Cells.Item(80, 8).Value = 2128.3157
$ws.Cells.Item(80, 9).Value = 629.2
$ws.Cells.Item(80, 10).Value = 3794
$ws.Cells.Item(80, 11).Value = 1887.6
$ws.Cells.Item(80, 12).Value = 11382
$ws.Cells.Item(80, 13).Value = -889.6000000000001
$ws.Cells.Item(80, 14).Value = -13378
$ws.Cells.Item(83, 8).Value = 2128.3157
$ws.Cells.Item(83, 9).Value = 629.2
$ws.Cells.Item(83, 10).Value = 3794
$ws.Cells.Item(83, 11).Value = 5662.8
$ws.Cells.Item(83, 12).Value = 34146
$ws.Cells.Item(83, 13).Value = -670.8000000000002
$ws.Cells.Item(83, 14).Value = -44130
$ws.Cells.Item(86, 8).Value = 9759.666999999999
$ws.Cells.Item(86, 9).Value = 9043.888999999999
$ws.Cells.Item(86, 10).Value = 10833.333
$ws.Cells.Item(86, 11).Value = 9043.888999999999
$ws.Cells.Item(86, 12).Value = 10833.333
$ws.Cells.Item(86, 13).Value = -7920.888999999999
$ws.Cells.Item(86, 14).Value = -13079.333
$ws.Cells.Item(89, 8).Value = 9759.666999999999
$ws.Cells.Item(89, 9).Value = 9043.888999999999
$ws.Cells.Item(89, 10).Value = 10833.333
$ws.Cells.Item(89, 11).Value = 45219.44499999999
$ws.Cells.Item(89, 12).Value = 54166.665
$ws.Cells.Item(89, 13).Value = -39603.44499999999
$ws.Cells.Item(89, 14).Value = -65398.665
$ws.Cells.Item(98, 8).Value = 1493.5217
$ws.Cells.Item(98, 9).Value = 1421.5238
$ws.Cells.Item(98, 11).Value = 1421.5238
$ws.Cells.Item(98, 13).Value = 76.47620000000006
$ws.Cells.Item(103, 8).Value = 879.0625
$ws.Cells.Item(103, 9).Value = 1026.5834
$ws.Cells.Item(103, 10).Value = 436.5
$ws.Cells.Item(103, 11).Value = 3079.7502
$ws.Cells.Item(103, 12).Value = 1309.5
$ws.Cells.Item(103, 13).Value = -2493.7502
$ws.Cells.Item(103, 14).Value = -2481.5
$ws.Cells.Item(106, 8).Value = 38470964
$ws.Cells.Item(106, 9).Value = 62512944
$ws.Cells.Item(106, 11).Value = 62512944
$ws.Cells.Item(106, 13).Value = -62512313
$ws.Cells.Item(109, 8).Value = 54000
$ws.Cells.Item(109, 10).Value = 0
$ws.Cells.Item(109, 12).Value = 0
$ws.Cells.Item(109, 14).ClearContents()
$ws.Cells.Item(122, 8).Value = 1493.5217
$ws.Cells.Item(122, 9).Value = 1421.5238
$ws.Cells.Item(122, 11).Value = 4264.5714
$ws.Cells.Item(122, 13).Value = -1814.5714
$ws.Cells.Item(131, 8).Value = 5166.724
$ws.Cells.Item(131, 9).Value = 2004.125
$ws.Cells.Item(131, 10).Value = 9059.154
$ws.Cells.Item(131, 11).Value = 6012.375
$ws.Cells.Item(131, 12).Value = 27177.462
$ws.Cells.Item(131, 13).Value = -972.375
$ws.Cells.Item(131, 14).Value = -37257.462
$ws.Cells.Item(137, 8).Value = 76785.03999999999
$ws.Cells.Item(137, 9).Value = 113366.75
$ws.Cells.Item(137, 10).Value = 3621.625
$ws.Cells.Item(137, 11).Value = 340100.25
$ws.Cells.Item(137, 12).Value = 10864.875
$ws.Cells.Item(137, 13).Value = -337550.25
$ws.Cells.Item(137, 14).Value = -15964.875

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(45, 8).Value = 5497663
$ws.Cells.Item(45, 9).Value = 10990079
$ws.Cells.Item(45, 11).Value = 10990079
$ws.Cells.Item(45, 13).Value = -10989702
$ws.Cells.Item(61, 8).Value = 6892.9414
$ws.Cells.Item(61, 9).Value = 7714.846
$ws.Cells.Item(61, 10).Value = 4221.75
$ws.Cells.Item(61, 11).Value = 7714.846
$ws.Cells.Item(61, 12).Value = 4221.75
$ws.Cells.Item(61, 13).Value = -7502.846
$ws.Cells.Item(61, 14).Value = -4645.75
$ws.Cells.Item(74, 8).Value = 43581.73
$ws.Cells.Item(74, 9).Value = 10667.235
$ws.Cells.Item(74, 10).Value = 105753.555
$ws.Cells.Item(74, 11).Value = 10667.235
$ws.Cells.Item(74, 12).Value = 105753.555
$ws.Cells.Item(74, 13).Value = -9793.235000000001
$ws.Cells.Item(74, 14).Value = -107501.555
$ws.Cells.Item(77, 8).Value = 43581.73
$ws.Cells.Item(77, 9).Value = 10667.235
$ws.Cells.Item(77, 10).Value = 105753.555
$ws.Cells.Item(77, 11).Value = 53336.175
$ws.Cells.Item(77, 12).Value = 528767.7749999999
$ws.Cells.Item(77, 13).Value = -48968.175
$ws.Cells.Item(77, 14).Value = -537503.7749999999
$ws.Cells.Item(122, 8).Value = 3502660.2
$ws.Cells.Item(122, 9).Value = 5265556.5
$ws.Cells.Item(122, 10).Value = 1900027.1
$ws.Cells.Item(122, 11).Value = 15796669.5
$ws.Cells.Item(122, 12).Value = 5700081.300000001
$ws.Cells.Item(122, 13).Value = -15794219.5
$ws.Cells.Item(122, 14).Value = -5704981.300000001
$ws.Cells.Item(136, 8).Value = 6892.9414
$ws.Cells.Item(136, 9).Value = 7714.846
$ws.Cells.Item(136, 10).Value = 4221.75
$ws.Cells.Item(136, 11).Value = 23144.538
$ws.Cells.Item(136, 12).Value = 12665.25
$ws.Cells.Item(136, 13).Value = -20594.538
$ws.Cells.Item(136, 14).Value = -17765.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 17340718
$ws.Cells.Item(86, 9).Value = 39396350
$ws.Cells.Item(86, 10).Value = 11290.143
$ws.Cells.Item(86, 11).Value = 39396350
$ws.Cells.Item(86, 12).Value = 11290.143
$ws.Cells.Item(86, 13).Value = -39395227
$ws.Cells.Item(86, 14).Value = -13536.143
$ws.Cells.Item(89, 8).Value = 17340718
$ws.Cells.Item(89, 9).Value = 39396350
$ws.Cells.Item(89, 10).Value = 11290.143
$ws.Cells.Item(89, 11).Value = 196981750
$ws.Cells.Item(89, 12).Value = 56450.715
$ws.Cells.Item(89, 13).Value = -196976134
$ws.Cells.Item(89, 14).Value = -67682.715

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(22, 8).Value = 924.3
$ws.Cells.Item(22, 9).Value = 1131.5333
$ws.Cells.Item(22, 10).Value = 302.6
$ws.Cells.Item(22, 11).Value = 1131.5333
$ws.Cells.Item(22, 12).Value = 302.6
$ws.Cells.Item(22, 13).Value = -781.5333000000001
$ws.Cells.Item(22, 14).Value = -1002.6
$ws.Cells.Item(31, 8).Value = 20549.334
$ws.Cells.Item(31, 9).Value = 5937.3184
$ws.Cells.Item(31, 10).Value = 29008.922
$ws.Cells.Item(31, 11).Value = 5937.3184
$ws.Cells.Item(31, 12).Value = 29008.922
$ws.Cells.Item(31, 13).Value = -5642.3184
$ws.Cells.Item(31, 14).Value = -29598.922
$ws.Cells.Item(34, 8).Value = 20549.334
$ws.Cells.Item(34, 9).Value = 5937.3184
$ws.Cells.Item(34, 10).Value = 29008.922
$ws.Cells.Item(34, 11).Value = 5937.3184
$ws.Cells.Item(34, 12).Value = 29008.922
$ws.Cells.Item(34, 13).Value = -5735.3184
$ws.Cells.Item(34, 14).Value = -29412.922
$ws.Cells.Item(99, 8).Value = 3751.389
$ws.Cells.Item(99, 10).Value = 4493.6665
$ws.Cells.Item(99, 12).Value = 4493.6665
$ws.Cells.Item(99, 14).Value = -7489.6665
$ws.Cells.Item(126, 8).Value = 3751.389
$ws.Cells.Item(126, 10).Value = 4493.6665
$ws.Cells.Item(126, 12).Value = 13480.9995
$ws.Cells.Item(126, 14).Value = -18420.9995
$ws.Cells.Item(132, 8).Value = 75415.19500000001
$ws.Cells.Item(132, 9).Value = 64628.375
$ws.Cells.Item(132, 10).Value = 92674.10000000001
$ws.Cells.Item(132, 11).Value = 193885.125
$ws.Cells.Item(132, 12).Value = 278022.3
$ws.Cells.Item(132, 13).Value = -191355.125
$ws.Cells.Item(132, 14).Value = -283082.3
$ws.Cells.Item(141, 8).Value = 148471.95
$ws.Cells.Item(141, 10).Value = 159122.1
$ws.Cells.Item(141, 12).Value = 159122.1
$ws.Cells.Item(141, 14).Value = -169482.1

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(33, 8).Value = 8447.833000000001
$ws.Cells.Item(33, 9).Value = 129.625
$ws.Cells.Item(33, 11).Value = 777.75
$ws.Cells.Item(33, 13).Value = -494.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(109, 8).Value = 0
$ws.Cells.Item(109, 10).Value = 0
$ws.Cells.Item(109, 12).Value = 0
$ws.Cells.Item(109, 14).ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 45620.95
$ws.Cells.Item(22, 9).Value = 68956.46000000001
$ws.Cells.Item(22, 10).Value = 2283.5715
$ws.Cells.Item(22, 11).Value = 68956.46000000001
$ws.Cells.Item(22, 12).Value = 2283.5715
$ws.Cells.Item(22, 13).Value = -68661.46000000001
$ws.Cells.Item(22, 14).Value = -2873.5715
$ws.Cells.Item(27, 8).Value = 45620.95
$ws.Cells.Item(27, 9).Value = 68956.46000000001
$ws.Cells.Item(27, 10).Value = 2283.5715
$ws.Cells.Item(27, 11).Value = 68956.46000000001
$ws.Cells.Item(27, 12).Value = 2283.5715
$ws.Cells.Item(27, 13).Value = -68849.46000000001
$ws.Cells.Item(27, 14).Value = -2497.5715
$ws.Cells.Item(46, 8).Value = 5679.375
$ws.Cells.Item(46, 9).Value = 449.5
$ws.Cells.Item(46, 10).Value = 6426.5
$ws.Cells.Item(46, 11).Value = 449.5
$ws.Cells.Item(46, 12).Value = 6426.5
$ws.Cells.Item(46, 13).Value = -261.5
$ws.Cells.Item(46, 14).Value = -6802.5
$ws.Cells.Item(55, 8).Value = 1503.3889
$ws.Cells.Item(55, 10).Value = 2317
$ws.Cells.Item(55, 12).Value = 2317
$ws.Cells.Item(55, 14).Value = -2663
$ws.Cells.Item(93, 8).Value = 9531659
$ws.Cells.Item(93, 9).Value = 14494874
$ws.Cells.Item(93, 10).Value = 18830.334
$ws.Cells.Item(93, 11).Value = 14494874
$ws.Cells.Item(93, 12).Value = 18830.334
$ws.Cells.Item(93, 13).Value = -14493626
$ws.Cells.Item(93, 14).Value = -21326.334
$ws.Cells.Item(100, 8).Value = 3608.4783
$ws.Cells.Item(100, 9).Value = 2999.4443
$ws.Cells.Item(100, 11).Value = 2999.4443
$ws.Cells.Item(100, 13).Value = -2458.4443
$ws.Cells.Item(132, 8).Value = 33032.547
$ws.Cells.Item(132, 9).Value = 35736.8
$ws.Cells.Item(132, 11).Value = 107210.4
$ws.Cells.Item(132, 13).Value = -104680.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(112, 8).Value = 35397.332
$ws.Cells.Item(112, 10).Value = 35397.332
$ws.Cells.Item(112, 12).Value = 35397.332
$ws.Cells.Item(112, 14).Value = -38351.332
$ws.Cells.Item(132, 8).Value = 14244900
$ws.Cells.Item(132, 9).Value = 15879151
$ws.Cells.Item(132, 10).Value = 1375170.5
$ws.Cells.Item(132, 11).Value = 47637453
$ws.Cells.Item(132, 12).Value = 4125511.5
$ws.Cells.Item(132, 13).Value = -47634923
$ws.Cells.Item(132, 14).Value = -4130571.5
$ws.Cells.Item(136, 8).Value = 6322.0938
$ws.Cells.Item(136, 9).Value = 6326.9
$ws.Cells.Item(136, 11).Value = 18980.7
$ws.Cells.Item(136, 13).Value = -16430.7
